$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 24 (this pushes the old rows 24-27 down to 25-28),
# picking up the formatting of row 23 above it.
$ws.Rows.Item(24).Insert()

# The site merge A19:A23 must now grow to cover the freshly inserted row 24.
$ws.Range("A19:A23").UnMerge()
$ws.Range("A19:A24").Merge()

# Rows 20-23 (and now 24) of this site block switch from the non-wrapping
# fill style to the wrapping fill style used at the top of the block.
$ws.Range("A20:A24").Style = "Normal"
$ws.Range("A20:A24").Interior.ColorIndex = $ws.Range("A19").Interior.ColorIndex
$ws.Range("A20:A24").HorizontalAlignment = $ws.Range("A19").HorizontalAlignment
$ws.Range("A20:A24").VerticalAlignment = $ws.Range("A19").VerticalAlignment
$ws.Range("A20:A24").WrapText = $true

# Fill in the new 2016 plot-setup record.
$ws.Range("B24").Value = 2016
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "2016 Apr"
$ws.Range("E24").Value = "spring wheat"
$ws.Range("F24").Value = "Keith (field manager?) on 2016-04-20"

$ws.Range("F25").Select()
